$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the BRR (C4) and SCI_BAUD (C5) values
$ws.Range("C4").Value = 162
$ws.Range("C5").Value = 9600

# Add extra baud-rate columns (D, E) for the SCI interrupt/BRR table
$ws.Range("D10").Value = 57600
$ws.Range("E10").Value = 9600

$ws.Range("D11").Formula = "=C8/((D10+1)*8)"
$ws.Range("E11").Formula = "=C8/((E10+1)*8)"

# Update row labels (set B11 before B10 so the shared-string table order matches)
$ws.Range("B11").Value = "?SCI_BRR?"
$ws.Range("B10").Value = "SCI BAUD"

# Update the active selection to the last filled cell
$ws.Range("E11").Select()
